# Generate Report for Handoff
# Updates the localization-status workbook for a new CI handoff run:
#  - the tracked .md file gets a fresh GUID name / handoff artifact names
#  - two new .png dependency files show up as "IsDependency" rows
#  - handoff timestamps are refreshed
#
# Overview   -> sheet 1 (file-level summary)
# zh-cn      -> sheet 2 (per-locale detail, zh-cn)
# de-de      -> sheet 3 (per-locale detail, de-de)

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# ---- new / changed identifiers -------------------------------------------------
$newMdGuid = "90463d6e-e0e8-4f2e-92d2-2c025e1246db"

$mdName   = "$newMdGuid.md"
$png1Name = "e7dec773-7448-4a29-a215-50b956a8b19a.png"
$png2Name = "ff7e7932-2f5e-48b9-8ab0-449510740b1e.png"

$zhXlfName = "$newMdGuid.1ac2fd46409a5907e8504ce0e3d344d5ba23d77c.zh-cn.xlf"
$deXlfName = "$newMdGuid.1ac2fd46409a5907e8504ce0e3d344d5ba23d77c.de-de.xlf"

$png1TargetName = "4a0d64d5692876053921fd58c0face122d0d3c2b.png"
$png2TargetName = "bd30aae3ace1a65702fb8e5c76aaff7ae51617a0.png"

$overviewDate = "2016-50-20 20:50:29"
$zhDate       = "2016-03-20 20:50:26"
$deDate       = "2016-03-20 20:50:29"
$epochDate    = "0001-01-01 00:00:00"

$statusReady        = "Ready for handoff"
$statusInclude      = "Include"
$statusIsDependency = "IsDependency"
$dependencyFrom     = "e2e\$mdName"

$mdHyperlink   = "https://github.com/OpenLocalizationTest/oltest/blob/c17993604db7aea3b1c9d01fcd49156f9cecc5d8/e2e/$mdName"
$png1Hyperlink = "https://github.com/OpenLocalizationTest/oltest/blob/c17993604db7aea3b1c9d01fcd49156f9cecc5d8/e2e/$png1Name"
$png2Hyperlink = "https://github.com/OpenLocalizationTest/oltest/blob/c17993604db7aea3b1c9d01fcd49156f9cecc5d8/e2e/$png2Name"

$zhXlfHyperlink = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/d5dfaa627263a18b66a31abe5547a76f0948837d/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$zhXlfName"
$deXlfHyperlink = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/dc4f2a07d6c3de18984926e7b5e057f6b67c9ff9/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$deXlfName"

$png1ZhTargetHyperlink = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/d5dfaa627263a18b66a31abe5547a76f0948837d/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$png1TargetName"
$png2ZhTargetHyperlink = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/d5dfaa627263a18b66a31abe5547a76f0948837d/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$png2TargetName"
$png1DeTargetHyperlink = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/dc4f2a07d6c3de18984926e7b5e057f6b67c9ff9/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$png1TargetName"
$png2DeTargetHyperlink = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/dc4f2a07d6c3de18984926e7b5e057f6b67c9ff9/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$png2TargetName"

$linkColor = 15570276   # RGB(100,149,237) == FF6495ED, matches the workbook's existing HyperLink style

function Style-LikeLink($range) {
    $range.Font.Underline = $true
    $range.Font.Color = $linkColor
}

function Style-LikeDatetime($range) {
    $range.NumberFormat = "yyyy-mm-dd HH:mm:ss"
}

# The engine's Hyperlink.Delete()/Address=/TextToDisplay= setters don't mutate
# an existing hyperlink record in place -- they leave the old one behind and
# append a second, duplicate entry. The collection-level Delete() does work,
# so each sheet's hyperlinks are cleared up front and rebuilt from scratch in
# final left-to-right, top-to-bottom order (this keeps rId1 == the table
# relationship, and rId2.. reassigned cleanly to the hyperlinks below).
$wsOverview.Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Delete()

# =====================================================================
# Overview sheet
# =====================================================================

$wsOverview.Range("A2").Value = $mdName
$wsOverview.Range("D2").Value = $overviewDate

$wsOverview.Range("A3").Value = $png1Name
$wsOverview.Range("B3").Value = $statusReady
$wsOverview.Range("C3").Value = $statusReady
$wsOverview.Range("D3").Value = $overviewDate

$wsOverview.Range("A4").Value = $png2Name
$wsOverview.Range("B4").Value = $statusReady
$wsOverview.Range("C4").Value = $statusReady
$wsOverview.Range("D4").Value = $overviewDate

$wsOverview.Hyperlinks.Add($wsOverview.Range("A2"), $mdHyperlink, [Type]::Missing, [Type]::Missing, $mdName) | Out-Null
$wsOverview.Hyperlinks.Add($wsOverview.Range("A3"), $png1Hyperlink, [Type]::Missing, [Type]::Missing, $png1Name) | Out-Null
$wsOverview.Hyperlinks.Add($wsOverview.Range("A4"), $png2Hyperlink, [Type]::Missing, [Type]::Missing, $png2Name) | Out-Null

Style-LikeLink($wsOverview.Range("A2"))
Style-LikeLink($wsOverview.Range("A3"))
Style-LikeLink($wsOverview.Range("A4"))

# =====================================================================
# zh-cn sheet
# =====================================================================

$wsZhCn.Range("A2").Value = $mdName
$wsZhCn.Range("D2").Value = $zhXlfName
$wsZhCn.Range("E2").Value = $zhDate

$wsZhCn.Range("A3").Value = $png1Name
$wsZhCn.Range("B3").Value = ".png"
$wsZhCn.Range("C3").Value = $statusReady
$wsZhCn.Range("D3").Value = $png1TargetName
$wsZhCn.Range("E3").Value = $zhDate
$wsZhCn.Range("H3").Value = $epochDate
$wsZhCn.Range("I3").Value = $statusIsDependency
$wsZhCn.Range("J3").Value = $dependencyFrom

$wsZhCn.Range("A4").Value = $png2Name
$wsZhCn.Range("B4").Value = ".png"
$wsZhCn.Range("C4").Value = $statusReady
$wsZhCn.Range("D4").Value = $png2TargetName
$wsZhCn.Range("E4").Value = $zhDate
$wsZhCn.Range("H4").Value = $epochDate
$wsZhCn.Range("I4").Value = $statusIsDependency
$wsZhCn.Range("J4").Value = $dependencyFrom

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), $mdHyperlink, [Type]::Missing, [Type]::Missing, $mdName) | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("B2"), $mdHyperlink, [Type]::Missing, [Type]::Missing, ".md") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("D2"), $zhXlfHyperlink, [Type]::Missing, [Type]::Missing, $zhXlfName) | Out-Null

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), $png1Hyperlink, [Type]::Missing, [Type]::Missing, $png1Name) | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("B3"), $png1Hyperlink, [Type]::Missing, [Type]::Missing, ".png") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("D3"), $png1ZhTargetHyperlink, [Type]::Missing, [Type]::Missing, $png1TargetName) | Out-Null

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A4"), $png2Hyperlink, [Type]::Missing, [Type]::Missing, $png2Name) | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("B4"), $png2Hyperlink, [Type]::Missing, [Type]::Missing, ".png") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("D4"), $png2ZhTargetHyperlink, [Type]::Missing, [Type]::Missing, $png2TargetName) | Out-Null

Style-LikeLink($wsZhCn.Range("A2"))
Style-LikeLink($wsZhCn.Range("B2"))
Style-LikeLink($wsZhCn.Range("D2"))
Style-LikeLink($wsZhCn.Range("A3"))
Style-LikeLink($wsZhCn.Range("B3"))
Style-LikeLink($wsZhCn.Range("D3"))
Style-LikeLink($wsZhCn.Range("A4"))
Style-LikeLink($wsZhCn.Range("B4"))
Style-LikeLink($wsZhCn.Range("D4"))

# =====================================================================
# de-de sheet
# =====================================================================

$wsDeDe.Range("A2").Value = $mdName
$wsDeDe.Range("D2").Value = $deXlfName
$wsDeDe.Range("E2").Value = $deDate

$wsDeDe.Range("A3").Value = $png1Name
$wsDeDe.Range("B3").Value = ".png"
$wsDeDe.Range("C3").Value = $statusReady
$wsDeDe.Range("D3").Value = $png1TargetName
$wsDeDe.Range("E3").Value = $deDate
$wsDeDe.Range("H3").Value = $epochDate
$wsDeDe.Range("I3").Value = $statusIsDependency
$wsDeDe.Range("J3").Value = $dependencyFrom

$wsDeDe.Range("A4").Value = $png2Name
$wsDeDe.Range("B4").Value = ".png"
$wsDeDe.Range("C4").Value = $statusReady
$wsDeDe.Range("D4").Value = $png2TargetName
$wsDeDe.Range("E4").Value = $deDate
$wsDeDe.Range("H4").Value = $epochDate
$wsDeDe.Range("I4").Value = $statusIsDependency
$wsDeDe.Range("J4").Value = $dependencyFrom

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), $mdHyperlink, [Type]::Missing, [Type]::Missing, $mdName) | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("B2"), $mdHyperlink, [Type]::Missing, [Type]::Missing, ".md") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("D2"), $deXlfHyperlink, [Type]::Missing, [Type]::Missing, $deXlfName) | Out-Null

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), $png1Hyperlink, [Type]::Missing, [Type]::Missing, $png1Name) | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("B3"), $png1Hyperlink, [Type]::Missing, [Type]::Missing, ".png") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("D3"), $png1DeTargetHyperlink, [Type]::Missing, [Type]::Missing, $png1TargetName) | Out-Null

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A4"), $png2Hyperlink, [Type]::Missing, [Type]::Missing, $png2Name) | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("B4"), $png2Hyperlink, [Type]::Missing, [Type]::Missing, ".png") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("D4"), $png2DeTargetHyperlink, [Type]::Missing, [Type]::Missing, $png2TargetName) | Out-Null

Style-LikeLink($wsDeDe.Range("A2"))
Style-LikeLink($wsDeDe.Range("B2"))
Style-LikeLink($wsDeDe.Range("D2"))
Style-LikeLink($wsDeDe.Range("A3"))
Style-LikeLink($wsDeDe.Range("B3"))
Style-LikeLink($wsDeDe.Range("D3"))
Style-LikeLink($wsDeDe.Range("A4"))
Style-LikeLink($wsDeDe.Range("B4"))
Style-LikeLink($wsDeDe.Range("D4"))
